$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 467.13208
$ws.Range("J17").Value = 467.13208
$ws.Range("L17").Value = 1401.39624
$ws.Range("N17").Value = -1737.39624
$ws.Range("H20").Value = 12125
$ws.Range("I20").Value = 12125
$ws.Range("K20").Value = 12125
$ws.Range("M20").Value = -11895
$ws.Range("H28").Value = 202.42857
$ws.Range("I28").Value = 198.5
$ws.Range("J28").Value = 210.28572
$ws.Range("K28").Value = 198.5
$ws.Range("L28").Value = 210.28572
$ws.Range("M28").Value = 286.5
$ws.Range("N28").Value = -1180.28572
$ws.Range("H33").Value = 241.5
$ws.Range("I33").Value = 248.6
$ws.Range("K33").Value = 248.6
$ws.Range("M33").Value = -19.59999999999999
$ws.Range("H35").Value = 12125
$ws.Range("I35").Value = 12125
$ws.Range("K35").Value = 12125
$ws.Range("M35").Value = -11746
$ws.Range("H70").Value = 2092.0833
$ws.Range("I70").Value = 3426.25
$ws.Range("J70").Value = 1425
$ws.Range("K70").Value = 10278.75
$ws.Range("L70").Value = 4275
$ws.Range("M70").Value = -10008.75
$ws.Range("N70").Value = -4815
$ws.Range("H73").Value = 2092.0833
$ws.Range("I73").Value = 3426.25
$ws.Range("J73").Value = 1425
$ws.Range("K73").Value = 10278.75
$ws.Range("L73").Value = 4275
$ws.Range("M73").Value = -9342.75
$ws.Range("N73").Value = -6147
$ws.Range("H96").Value = 1107.091
$ws.Range("I96").Value = 350
$ws.Range("J96").Value = 1738
$ws.Range("K96").Value = 1050
$ws.Range("L96").Value = 5214
$ws.Range("M96").Value = 323
$ws.Range("N96").Value = -7960
$ws.Range("H112").Value = 980.1896400000001
$ws.Range("J112").Value = 997.3393
$ws.Range("L112").Value = 2992.0179
$ws.Range("N112").Value = -5208.0179
$ws.Range("H127").Value = 1346.9166
$ws.Range("I127").Value = 297
$ws.Range("J127").Value = 2396.8333
$ws.Range("K127").Value = 891
$ws.Range("L127").Value = 7190.499899999999
$ws.Range("M127").Value = 4069
$ws.Range("N127").Value = -17110.4999
$ws.Range("H138").Value = 1514.6198
$ws.Range("I138").Value = 607.413
$ws.Range("J138").Value = 3183.88
$ws.Range("K138").Value = 1822.239
$ws.Range("L138").Value = 9551.639999999999
$ws.Range("M138").Value = 3317.761
$ws.Range("N138").Value = -19831.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 921.6111
$ws.Range("I2").Value = 926
$ws.Range("J2").Value = 914.7143
$ws.Range("K2").Value = 926
$ws.Range("L2").Value = 914.7143
$ws.Range("M2").Value = -813
$ws.Range("N2").Value = -1140.7143
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H116").Value = 921.6111
$ws.Range("I116").Value = 926
$ws.Range("J116").Value = 914.7143
$ws.Range("K116").Value = 926
$ws.Range("L116").Value = 914.7143
$ws.Range("M116").Value = 1368
$ws.Range("N116").Value = -5502.7143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 921.6111
$ws.Range("I3").Value = 926
$ws.Range("J3").Value = 914.7143
$ws.Range("K3").Value = 926
$ws.Range("L3").Value = 914.7143
$ws.Range("M3").Value = -812
$ws.Range("N3").Value = -1142.7143
$ws.Range("H86").Value = 970460.9399999999
$ws.Range("I86").Value = 1279.4706
$ws.Range("J86").Value = 3324187.2
$ws.Range("K86").Value = 1279.4706
$ws.Range("L86").Value = 3324187.2
$ws.Range("M86").Value = -156.4706000000001
$ws.Range("N86").Value = -3326433.2
$ws.Range("H89").Value = 970460.9399999999
$ws.Range("I89").Value = 1279.4706
$ws.Range("J89").Value = 3324187.2
$ws.Range("K89").Value = 6397.353000000001
$ws.Range("L89").Value = 16620936
$ws.Range("M89").Value = -781.353000000001
$ws.Range("N89").Value = -16632168
$ws.Range("H94").Value = 939.6
$ws.Range("I94").Value = 924.5
$ws.Range("K94").Value = 924.5
$ws.Range("M94").Value = -473.5
$ws.Range("H99").Value = 2199.5
$ws.Range("J99").Value = 2472.0908
$ws.Range("L99").Value = 2472.0908
$ws.Range("N99").Value = -5468.0908
$ws.Range("H134").Value = 2060721.8
$ws.Range("I134").Value = 1138.1316
$ws.Range("J134").Value = 6952233
$ws.Range("K134").Value = 3414.3948
$ws.Range("L134").Value = 20856699
$ws.Range("M134").Value = -879.3948
$ws.Range("N134").Value = -20861769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 884428.25
$ws.Range("I31").Value = 1012732.1
$ws.Range("J31").Value = 2339.125
$ws.Range("K31").Value = 1012732.1
$ws.Range("L31").Value = 2339.125
$ws.Range("M31").Value = -1012437.1
$ws.Range("N31").Value = -2929.125
$ws.Range("H34").Value = 884428.25
$ws.Range("I34").Value = 1012732.1
$ws.Range("J34").Value = 2339.125
$ws.Range("K34").Value = 1012732.1
$ws.Range("L34").Value = 2339.125
$ws.Range("M34").Value = -1012530.1
$ws.Range("N34").Value = -2743.125
$ws.Range("H99").Value = 3530.5715
$ws.Range("I99").Value = 2557.1428
$ws.Range("J99").Value = 4504
$ws.Range("K99").Value = 2557.1428
$ws.Range("L99").Value = 4504
$ws.Range("M99").Value = -1059.1428
$ws.Range("N99").Value = -7500
$ws.Range("H122").Value = 31250626
$ws.Range("I122").Value = 35714784
$ws.Range("K122").Value = 107144352
$ws.Range("M122").Value = -107141902
$ws.Range("H126").Value = 3530.5715
$ws.Range("I126").Value = 2557.1428
$ws.Range("J126").Value = 4504
$ws.Range("K126").Value = 7671.428400000001
$ws.Range("L126").Value = 13512
$ws.Range("M126").Value = -5201.428400000001
$ws.Range("N126").Value = -18452
$ws.Range("H132").Value = 27779888
$ws.Range("I132").Value = 1479
$ws.Range("J132").Value = 111115110
$ws.Range("K132").Value = 4437
$ws.Range("L132").Value = 333345330
$ws.Range("M132").Value = -1907
$ws.Range("N132").Value = -333350390

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 253
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 253
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 759
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1229
$ws.Range("H114").Value = 613.7143
$ws.Range("I114").Value = 206.33333
$ws.Range("J114").Value = 776.6667
$ws.Range("K114").Value = 618.99999
$ws.Range("L114").Value = 2330.0001
$ws.Range("M114").Value = 2635.00001
$ws.Range("N114").Value = -8838.000100000001
$ws.Range("H121").Value = 111111110
$ws.Range("J121").Value = 111111110
$ws.Range("L121").Value = 333333330
$ws.Range("N121").Value = -333335950
$ws.Range("H134").Value = 23812350
$ws.Range("I134").Value = 45455396
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 136366188
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = -136361118
$ws.Range("N134").Value = -25137

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 12250000
$ws.Range("I11").Value = 8000000
$ws.Range("J11").Value = 25000000
$ws.Range("K11").Value = 8000000
$ws.Range("L11").Value = 25000000
$ws.Range("M11").Value = -7999861
$ws.Range("N11").Value = -25000278
$ws.Range("H12").Value = 3502163.2
$ws.Range("I12").Value = 4201596
$ws.Range("J12").Value = 5000
$ws.Range("K12").Value = 4201596
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = -4201456
$ws.Range("N12").Value = -5280
$ws.Range("H14").Value = 173.57143
$ws.Range("I14").Value = 173.57143
$ws.Range("K14").Value = 173.57143
$ws.Range("M14").Value = -5.571429999999992
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H80").Value = 33334734
$ws.Range("I80").Value = 2100
$ws.Range("K80").Value = 2100
$ws.Range("M80").Value = -1102
$ws.Range("H83").Value = 33334734
$ws.Range("I83").Value = 2100
$ws.Range("K83").Value = 10500
$ws.Range("M83").Value = -5508
$ws.Range("H102").Value = 1038.5555
$ws.Range("I102").Value = 1069.409
$ws.Range("J102").Value = 902.8
$ws.Range("K102").Value = 1069.409
$ws.Range("L102").Value = 902.8
$ws.Range("M102").Value = 552.5909999999999
$ws.Range("N102").Value = -4146.8
$ws.Range("H126").Value = 2323.2
$ws.Range("I126").Value = 1340
$ws.Range("J126").Value = 3446.8572
$ws.Range("K126").Value = 4020
$ws.Range("L126").Value = 10340.5716
$ws.Range("M126").Value = -1550
$ws.Range("N126").Value = -15280.5716
$ws.Range("H132").Value = 5514.032
$ws.Range("I132").Value = 2076.15
$ws.Range("J132").Value = 11764.728
$ws.Range("K132").Value = 6228.450000000001
$ws.Range("L132").Value = 35294.18399999999
$ws.Range("M132").Value = -3698.450000000001
$ws.Range("N132").Value = -40354.18399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 41209136
$ws.Range("I16").Value = 3571769.5
$ws.Range("J16").Value = 166667020
$ws.Range("K16").Value = 3571769.5
$ws.Range("L16").Value = 166667020
$ws.Range("M16").Value = -3571599.5
$ws.Range("N16").Value = -166667360
$ws.Range("H61").Value = 1427.9286
$ws.Range("I61").Value = 1429.3
$ws.Range("J61").Value = 1424.5
$ws.Range("K61").Value = 1429.3
$ws.Range("L61").Value = 1424.5
$ws.Range("M61").Value = -1227.3
$ws.Range("N61").Value = -1828.5
$ws.Range("H82").Value = 9524867
$ws.Range("I82").Value = 10205115
$ws.Range("K82").Value = 10205115
$ws.Range("M82").Value = -10204754
$ws.Range("H85").Value = 9524867
$ws.Range("I85").Value = 10205115
$ws.Range("K85").Value = 10205115
$ws.Range("M85").Value = -10203867
$ws.Range("H113").Value = 1427.9286
$ws.Range("I113").Value = 1429.3
$ws.Range("J113").Value = 1424.5
$ws.Range("K113").Value = 1429.3
$ws.Range("L113").Value = 1424.5
$ws.Range("M113").Value = 740.7
$ws.Range("N113").Value = -5764.5
$ws.Range("H136").Value = 42330620
$ws.Range("I136").Value = 5955197
$ws.Range("K136").Value = 17865591
$ws.Range("M136").Value = -17863041

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 21659.076
$ws.Range("I122").Value = 37181.855
$ws.Range("J122").Value = 3549.1667
$ws.Range("K122").Value = 111545.565
$ws.Range("L122").Value = 10647.5001
$ws.Range("M122").Value = -109095.565
$ws.Range("N122").Value = -15547.5001
$ws.Range("H126").Value = 726.9
$ws.Range("I126").Value = 455.17648
$ws.Range("J126").Value = 2266.6667
$ws.Range("K126").Value = 1365.52944
$ws.Range("L126").Value = 6800.000100000001
$ws.Range("M126").Value = 1104.47056
$ws.Range("N126").Value = -11740.0001
$ws.Range("H132").Value = 22487.453
$ws.Range("I132").Value = 28702.352
$ws.Range("J132").Value = 8115.5
$ws.Range("K132").Value = 86107.056
$ws.Range("L132").Value = 24346.5
$ws.Range("M132").Value = -83577.056
$ws.Range("N132").Value = -29406.5

